# Generate Report for Handoff
# Updates the status (and handoff timestamps) for the
# fa496268-6053-49df-8127-4f8d4c5800aa.md file to reflect that it is now
# ready for handoff, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the fa496268-...md entry (B=zh-cn status, C=de-de status)
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the fa496268-...md entry
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-28 05:28:57"

# de-de sheet: row 3 is the fa496268-...md entry
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-28 05:29:08"
